# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mirrors a "handback" report generation: the two language
# sheets (zh-cn / de-de) gain a "Latest Target File" (F) and "Latest
# Handback File" (G) hyperlinked entry for each row, the overall Status
# text flips from "Ready for handoff" to "Handed back: in sync with
# en-US" everywhere it appears, and the "Latest Handback DateTime" (H)
# column is stamped with real handback timestamps instead of the
# zero-date placeholder.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Status text: every cell that currently reads "Ready for handoff"
#    now reads "Handed back: in sync with en-US" (Overview!B/C2:3 and
#    the Status column (C) on both language sheets).
# ---------------------------------------------------------------------
$ws_overview.Range("B2").Value = $newStatus
$ws_overview.Range("C2").Value = $newStatus
$ws_overview.Range("B3").Value = $newStatus
$ws_overview.Range("C3").Value = $newStatus

$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus

$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) zh-cn sheet: add "Latest Target File" (F) / "Latest Handback File"
#    (G) hyperlinked values, and stamp the handback datetime (H).
# ---------------------------------------------------------------------
$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b0b4a79c82b0454d140d59ed0b45b3adbb8e9c6/e2e/7476a706-939e-4cf3-b8d0-b63bde0bbb88.md",
    "",
    "",
    "7476a706-939e-4cf3-b8d0-b63bde0bbb88.md") | Out-Null

$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c949c5f5e8b59d61f9174f77340999fd2624c59d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/7476a706-939e-4cf3-b8d0-b63bde0bbb88.bc5a12cf969f33e1611ce08c652f642de4ec15e6.zh-cn.xlf",
    "",
    "",
    "7476a706-939e-4cf3-b8d0-b63bde0bbb88.bc5a12cf969f33e1611ce08c652f642de4ec15e6.zh-cn.xlf") | Out-Null

$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b0b4a79c82b0454d140d59ed0b45b3adbb8e9c6/e2e/b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.md",
    "",
    "",
    "b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.md") | Out-Null

$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c949c5f5e8b59d61f9174f77340999fd2624c59d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.fd6a8a19b03b8c7a19d8d0ca6b01be594d074a39.zh-cn.xlf",
    "",
    "",
    "b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.fd6a8a19b03b8c7a19d8d0ca6b01be594d074a39.zh-cn.xlf") | Out-Null

$ws_zhcn.Range("H2").Value = "2016-03-12 22:12:03"
$ws_zhcn.Range("H3").Value = "2016-03-12 22:12:03"

# ---------------------------------------------------------------------
# 3) de-de sheet: same additions, de-de variants, with its own handback
#    timestamp.
# ---------------------------------------------------------------------
$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b0b4a79c82b0454d140d59ed0b45b3adbb8e9c6/e2e/7476a706-939e-4cf3-b8d0-b63bde0bbb88.md",
    "",
    "",
    "7476a706-939e-4cf3-b8d0-b63bde0bbb88.md") | Out-Null

$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e497ea16ba70e9505070b7ed7dbdaa12fc485f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/7476a706-939e-4cf3-b8d0-b63bde0bbb88.bc5a12cf969f33e1611ce08c652f642de4ec15e6.de-de.xlf",
    "",
    "",
    "7476a706-939e-4cf3-b8d0-b63bde0bbb88.bc5a12cf969f33e1611ce08c652f642de4ec15e6.de-de.xlf") | Out-Null

$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8b0b4a79c82b0454d140d59ed0b45b3adbb8e9c6/e2e/b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.md",
    "",
    "",
    "b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.md") | Out-Null

$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e497ea16ba70e9505070b7ed7dbdaa12fc485f2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.fd6a8a19b03b8c7a19d8d0ca6b01be594d074a39.de-de.xlf",
    "",
    "",
    "b8a959ff-78d7-4c2d-9ade-20a0d4d13acf.fd6a8a19b03b8c7a19d8d0ca6b01be594d074a39.de-de.xlf") | Out-Null

$ws_dede.Range("H2").Value = "2016-03-12 22:12:10"
$ws_dede.Range("H3").Value = "2016-03-12 22:12:10"
